$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.441.70"
$ws.Range("E2").Value = "  +4.01%  "
$ws.Range("D3").Value = "2.464.97"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +2.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.93"
$ws.Range("E6").Value = "  +4.57%  "
$ws.Range("E7").Value = "  +1.28%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("E9").Value = "  +3.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.30"
$ws.Range("E10").Value = "  +2.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0817"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.43"
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.10"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Value = "2.850.12"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").Value = "2.478.16"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.848"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "46.315.80"
$ws.Range("E18").Value = "  +4.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.76"
$ws.Range("E19").Value = "  +2.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.45"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").Value = "0.0₃0938"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.59"
$ws.Range("E22").Value = "  +2.71%  "
$ws.Range("E23").Value = "  +4.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "248.47"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("E25").Value = "  +2.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.18"
$ws.Range("E26").Value = "  +3.68%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  -2.72%  "
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.68"
$ws.Range("E30").Value = "  +3.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "49.57"
$ws.Range("E31").Value = "  +2.32%  "
$ws.Range("E32").Value = "  +3.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.83"
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.35"
$ws.Range("E34").Value = "  +3.19%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.91"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "123.13"
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.112"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.80"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("D45").Value = "1.980.54"
$ws.Range("E45").Value = "  +2.04%  "
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("E47").Value = "  -3.06%  "
$ws.Range("E48").Value = "  +10.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.11"
$ws.Range("E49").Value = "  -3.88%  "
$ws.Range("E50").Value = "  +10.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.77"
$ws.Range("E51").Value = "  +4.27%  "
